$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2825
$ws.Range("H64").Value = 4103.1665
$ws.Range("I64").Value = 4164.75
$ws.Range("K64").Value = 4164.75
$ws.Range("M64").Value = -3916.75
$ws.Range("H67").Value = 4103.1665
$ws.Range("I67").Value = 4164.75
$ws.Range("K67").Value = 4164.75
$ws.Range("M67").Value = -3306.75
$ws.Range("H76").Value = 5629.7
$ws.Range("I76").Value = 5374.75
$ws.Range("K76").Value = 5374.75
$ws.Range("M76").Value = -5059.75
$ws.Range("H79").Value = 5629.7
$ws.Range("I79").Value = 5374.75
$ws.Range("K79").Value = 5374.75
$ws.Range("M79").Value = -4282.75
$ws.Range("H103").Value = 1753.3334
$ws.Range("J103").Value = 1940
$ws.Range("L103").Value = 5820
$ws.Range("N103").Value = -6992
$ws.Range("H112").Value = 2478.6
$ws.Range("I112").Value = 699.25
$ws.Range("K112").Value = 2097.75
$ws.Range("M112").Value = -989.75
$ws.Range("H129").Value = 801.61536
$ws.Range("J129").Value = 873.51514
$ws.Range("L129").Value = 2620.54542
$ws.Range("N129").Value = -12620.54542
$ws.Range("H138").Value = 1978.2424
$ws.Range("I138").Value = 1492.9584
$ws.Range("J138").Value = 2133.5334
$ws.Range("K138").Value = 4478.8752
$ws.Range("L138").Value = 6400.600199999999
$ws.Range("M138").Value = 661.1247999999996
$ws.Range("N138").Value = -16680.6002

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 28572478
$ws.Range("I61").Value = 32258898
$ws.Range("K61").Value = 32258898
$ws.Range("M61").Value = -32258686
$ws.Range("H63").Value = 18870580
$ws.Range("I63").Value = 2292.4783
$ws.Range("K63").Value = 2292.4783
$ws.Range("M63").Value = -1606.4783
$ws.Range("H66").Value = 18870580
$ws.Range("I66").Value = 2292.4783
$ws.Range("K66").Value = 11462.3915
$ws.Range("M66").Value = -8030.391500000002
$ws.Range("H88").Value = 3175
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 3566.6667
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 3566.6667
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -4378.6667
$ws.Range("H91").Value = 3175
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 3566.6667
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 3566.6667
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -6374.6667
$ws.Range("H132").Value = 2289.5615
$ws.Range("I132").Value = 1648.5
$ws.Range("J132").Value = 4084.5334
$ws.Range("K132").Value = 4945.5
$ws.Range("L132").Value = 12253.6002
$ws.Range("M132").Value = -2415.5
$ws.Range("N132").Value = -17313.6002
$ws.Range("H136").Value = 28572478
$ws.Range("I136").Value = 32258898
$ws.Range("K136").Value = 96776694
$ws.Range("M136").Value = -96774144

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 83334760
$ws.Range("I105").Value = 100001110
$ws.Range("K105").Value = 100001110
$ws.Range("M105").Value = -99999363

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 15000
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15480
$ws.Range("H27").Value = 15000
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 15000
$ws.Range("N27").Value = -15384
$ws.Range("H62").Value = 13338247
$ws.Range("J62").Value = 40001840
$ws.Range("L62").Value = 40001840
$ws.Range("N62").Value = -40003088
$ws.Range("H65").Value = 13338247
$ws.Range("J65").Value = 40001840
$ws.Range("L65").Value = 200009200
$ws.Range("N65").Value = -200015440
$ws.Range("H99").Value = 1586.1904
$ws.Range("I99").Value = 1515.625
$ws.Range("J99").Value = 1812
$ws.Range("K99").Value = 1515.625
$ws.Range("L99").Value = 1812
$ws.Range("M99").Value = -17.625
$ws.Range("N99").Value = -4808
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H126").Value = 1586.1904
$ws.Range("I126").Value = 1515.625
$ws.Range("J126").Value = 1812
$ws.Range("K126").Value = 4546.875
$ws.Range("L126").Value = 5436
$ws.Range("M126").Value = -2076.875
$ws.Range("N126").Value = -10376
$ws.Range("N110").ClearContents()
$ws.Range("N111").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1266.6666
$ws.Range("I17").Value = 1336.3636
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 4009.0908
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -3840.0908
$ws.Range("N17").Value = -1838
$ws.Range("H39").Value = 3960.8
$ws.Range("J39").Value = 3985.5386
$ws.Range("L39").Value = 11956.6158
$ws.Range("N39").Value = -12544.6158
$ws.Range("H55").Value = 2406.5386
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2406.5386
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 7219.6158
$ws.Range("N55").Value = -7573.6158
$ws.Range("H113").Value = 668.9545000000001
$ws.Range("I113").Value = 582.53845
$ws.Range("J113").Value = 705.19354
$ws.Range("K113").Value = 1747.61535
$ws.Range("L113").Value = 2115.58062
$ws.Range("M113").Value = 422.38465
$ws.Range("N113").Value = -6455.58062
$ws.Range("H131").Value = 12840490
$ws.Range("I131").Value = 76923480
$ws.Range("J131").Value = 23892.615
$ws.Range("K131").Value = 230770440
$ws.Range("L131").Value = 71677.845
$ws.Range("M131").Value = -230765400
$ws.Range("N131").Value = -81757.845
$ws.Range("M55").ClearContents()

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30003718
$ws.Range("I70").Value = 20837246
$ws.Range("K70").Value = 20837246
$ws.Range("M70").Value = -20836976
$ws.Range("H73").Value = 30003718
$ws.Range("I73").Value = 20837246
$ws.Range("K73").Value = 20837246
$ws.Range("M73").Value = -20836310
$ws.Range("H80").Value = 6800
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6800
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 6800
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 34000
$ws.Range("N83").Value = -43984
$ws.Range("H105").Value = 26999.5
$ws.Range("J105").Value = 26999.5
$ws.Range("L105").Value = 26999.5
$ws.Range("N105").Value = -33987.5
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4849.0835
$ws.Range("I46").Value = 797.5
$ws.Range("J46").Value = 6874.875
$ws.Range("K46").Value = 797.5
$ws.Range("L46").Value = 6874.875
$ws.Range("M46").Value = -609.5
$ws.Range("N46").Value = -7250.875
